$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 80: pnorm / z-test row
$ws.Range("A80").Value = "分布"
$ws.Range("C80").Value = "pnorm(z, lower.tail = a)"
$ws.Range("D80").Value = "z: 检验统计量`na: 布尔值，当TRUE时为双边检测，当FALSE时为单边检测"
$ws.Range("D80").Characters(4,6).Font.Name = "宋体"
$ws.Range("D80").Characters(10,3).Font.Name = "Cascadia Code"
$ws.Range("D80").Characters(13,5).Font.Name = "宋体"
$ws.Range("D80").Characters(18,4).Font.Name = "Cascadia Code"
$ws.Range("D80").Characters(22,8).Font.Name = "宋体"
$ws.Range("D80").Characters(30,5).Font.Name = "Cascadia Code"
$ws.Range("D80").Characters(35,6).Font.Name = "宋体"
$ws.Range("E80").Value = "进行单样本p值假设检验（计算p值）"
$ws.Range("E80").Characters(1,5).Font.Name = "等线"
$ws.Range("E80").Characters(6,1).Font.Name = "Cascadia Code"
$ws.Range("E80").Characters(7,8).Font.Name = "等线"
$ws.Range("E80").Characters(15,1).Font.Name = "Cascadia Code"
$ws.Range("E80").Characters(16,2).Font.Name = "等线"
$ws.Rows.Item(80).RowHeight = 48

# Row 81: t.test row
$ws.Range("A81").Value = "分布"
$ws.Range("C81").Value = "t.test(x, mu = a, alternative = b)"
$ws.Range("D81").Value = "x: 一个向量，或数据表中的一列，或某个numeric对象`na: 预先指定值（零假设中正确的概率）`nb: 字符串，当为`"greater`"或`"less`"时为单边检测"
$ws.Range("D81").Characters(4,17).Font.Name = "等线"
$ws.Range("D81").Characters(21,7).Font.Name = "Cascadia Code"
$ws.Range("D81").Characters(28,2).Font.Name = "等线"
$ws.Range("D81").Characters(30,4).Font.Name = "Cascadia Code"
$ws.Range("D81").Characters(34,16).Font.Name = "宋体"
$ws.Range("D81").Characters(50,4).Font.Name = "Cascadia Code"
$ws.Range("D81").Characters(54,6).Font.Name = "宋体"
$ws.Range("D81").Characters(60,9).Font.Name = "Cascadia Code"
$ws.Range("D81").Characters(69,1).Font.Name = "宋体"
$ws.Range("D81").Characters(70,6).Font.Name = "Cascadia Code"
$ws.Range("D81").Characters(76,6).Font.Name = "宋体"
$ws.Range("E81").Value = "进行单样本t检验（计算t值）。在样本容`n量足够大（不少于100）的情况下，t检验和z检验的结果差距不大"
$ws.Range("E81").Characters(1,5).Font.Name = "等线"
$ws.Range("E81").Characters(6,1).Font.Name = "Cascadia Code"
$ws.Range("E81").Characters(7,5).Font.Name = "等线"
$ws.Range("E81").Characters(12,1).Font.Name = "Cascadia Code"
$ws.Range("E81").Characters(13,2).Font.Name = "等线"
$ws.Range("E81").Characters(15,5).Font.Name = "宋体"
$ws.Range("E81").Characters(20,1).Font.Name = "Cascadia Code"
$ws.Range("E81").Characters(21,8).Font.Name = "宋体"
$ws.Range("E81").Characters(29,3).Font.Name = "Cascadia Code"
$ws.Range("E81").Characters(32,6).Font.Name = "宋体"
$ws.Range("E81").Characters(38,1).Font.Name = "Cascadia Code"
$ws.Range("E81").Characters(39,3).Font.Name = "宋体"
$ws.Range("E81").Characters(42,1).Font.Name = "Cascadia Code"
$ws.Range("E81").Characters(43,9).Font.Name = "宋体"
$ws.Rows.Item(81).RowHeight = 78

# Update the active selection to match the end-state of the editing session
$ws.Range("E87").Select()
